# Applies the "Trade #35 closed" update to the live trading results workbook.
#
# Summary of changes:
#  - Summary sheet: update Current Capital, Total P&L $, Total P&L %,
#    Total Trades, Winning Trades, Win Rate %
#  - Strategy Status sheet: update MarketMaking row (Capital, Trades, P&L $,
#    P&L %, Win Rate %)
#  - All Trades & MarketMaking sheets: append new trade row (#35)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.56   # Current Capital
$summary.Range("B4").Value = -0.44     # Total P&L $
$summary.Range("B5").Value = -0.25     # Total P&L %
$summary.Range("B6").Value = 35        # Total Trades
$summary.Range("B7").Value = 10        # Winning Trades
$summary.Range("B9").Value = 28.57     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.56      # Capital
$status.Range("D4").Value = 35         # Trades
$status.Range("E4").Value = -0.44      # P&L $
$status.Range("F4").Value = -0.44      # P&L %
$status.Range("G4").Value = 28.57      # Win Rate %

# ---------------------------------------------------------------------
# Helper to append the new trade row (#35) to a trades log sheet
# ---------------------------------------------------------------------
function Add-TradeRow35($sheet) {
    $sheet.Range("A36").Value = 35

    # The Date column stores plain text like "2026-02-17" (matching the
    # rest of the column) rather than a real date - force text entry so
    # Excel does not auto-convert it to a date serial number, then
    # restore the default "Normal" style so no stray formatting is left
    # on the cell.
    $sheet.Range("B36").NumberFormat = "@"
    $sheet.Range("B36").Value = "2026-02-17"
    $sheet.Range("B36").Style = "Normal"

    $sheet.Range("C36").Value = "08:32:37"
    $sheet.Range("D36").Value = "MarketMaking"
    $sheet.Range("E36").Value = "UP"
    $sheet.Range("F36").Value = 0.29
    $sheet.Range("G36").Value = 0.39604
    $sheet.Range("H36").Value = "CLOSED"
    $sheet.Range("I36").Value = 36.5654
    $sheet.Range("J36").Value = 0.11
    $sheet.Range("K36").Value = 99.56
    $sheet.Range("L36").Value = 0
    $sheet.Range("M36").Value = 0
    $sheet.Range("N36").Value = 0.6
    $sheet.Range("O36").Value = "Normal spread capture: 19600 bps"
    $sheet.Range("P36").Value = "early_exit"
    $sheet.Range("Q36").Value = 0.13
}

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow35 $allTrades

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow35 $marketMaking
